# Weekly refresh of Fruta / Hortaliza prices: the rows' Fecha / Calidad /
# Volumen / Precio mínimo / Precio máximo / Precio promedio ponderado /
# Precio $/Kg columns are re-shuffled among the 18 data rows (2-19).
# Columns A,B,C,E,F,G,H,N,O,Q,R are unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row -> source row (values that row now receives
# are the values that used to live in the source row).
$mapping = @{
    2  = 10
    3  = 18
    4  = 19
    5  = 8
    6  = 11
    7  = 12
    8  = 13
    9  = 16
    10 = 9
    11 = 6
    12 = 7
    13 = 4
    14 = 5
    15 = 3
    16 = 14
    17 = 15
    18 = 17
    19 = 2
}

$cols = @("D", "I", "J", "K", "L", "M", "P")

# Snapshot the original values for the columns that move, before any
# writes happen, so overlapping writes don't clobber source data.
$snapshot = @{}
foreach ($row in 2..19) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $rowVals
}

foreach ($destRow in 2..19) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $srcVals[$col]
    }
}
